$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("M2").Value = 0.303146
$ws.Range("N2").Value = 0.909438
$ws.Range("O2").Value = 0.005142855213700541
$ws.Range("P2").Value = 0.005142855213700542
$ws.Range("Q2").Value = 0.16322744797
$ws.Range("R2").Value = 1.46904703173
$ws.Range("S2").Value = 0.0001733969179253629
$ws.Range("T2").Value = 0.0001733969179253629
$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("O3").Value = 0.2877784259203595
$ws.Range("P3").Value = 0.2877784259203595
$ws.Range("Q3").Value = 9.133708045808332
$ws.Range("R3").Value = 82.203372412275
$ws.Range("S3").Value = 0.009702760436861925
$ws.Range("T3").Value = 0.009702760436861925
$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 41.67881
$ws.Range("N4").Value = 125.03643
$ws.Range("O4").Value = 0.7070787188659401
$ws.Range("P4").Value = 0.7070787188659401
$ws.Range("Q4").Value = 22.44174685045
$ws.Range("R4").Value = 201.97572165405
$ws.Range("S4").Value = 0.02383992266695518
$ws.Range("T4").Value = 0.02383992266695518
$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("M5").Value = 0.303146
$ws.Range("N5").Value = 0.909438
$ws.Range("O5").Value = 0.005142855213700541
$ws.Range("P5").Value = 0.005142855213700542
$ws.Range("Q5").Value = 3.650008091731333
$ws.Range("R5").Value = 32.850072825582
$ws.Range("S5").Value = 0.003877412539251186
$ws.Range("T5").Value = 0.003877412539251186
$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("O6").Value = 0.2877784259203595
$ws.Range("P6").Value = 0.2877784259203595
$ws.Range("Q6").Value = 204.2432733546095
$ws.Range("S6").Value = 0.2169681297301525
$ws.Range("T6").Value = 0.2169681297301525
$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 41.67881
$ws.Range("N7").Value = 125.03643
$ws.Range("O7").Value = 0.7070787188659401
$ws.Range("P7").Value = 0.7070787188659401
$ws.Range("Q7").Value = 501.8307803953634
$ws.Range("R7").Value = 4516.47702355827
$ws.Range("S7").Value = 0.5330960676211057
$ws.Range("T7").Value = 0.5330960676211057
$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("M8").Value = 0.303146
$ws.Range("N8").Value = 0.909438
$ws.Range("O8").Value = 0.005142855213700541
$ws.Range("P8").Value = 0.005142855213700542
$ws.Range("Q8").Value = 1.027998905843333
$ws.Range("R8").Value = 9.251990152589999
$ws.Range("S8").Value = 0.001092045756523993
$ws.Range("T8").Value = 0.001092045756523993
$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("O9").Value = 0.2877784259203595
$ws.Range("P9").Value = 0.2877784259203595
$ws.Range("Q9").Value = 57.52367015570277
$ws.Range("R9").Value = 517.7130314013249
$ws.Range("S9").Value = 0.06110753575334504
$ws.Range("T9").Value = 0.06110753575334503
$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 41.67881
$ws.Range("N10").Value = 125.03643
$ws.Range("O10").Value = 0.7070787188659401
$ws.Range("P10").Value = 0.7070787188659401
$ws.Range("Q10").Value = 141.3370820556833
$ws.Range("R10").Value = 1272.03373850115
$ws.Range("S10").Value = 0.1501427285778792
$ws.Range("T10").Value = 0.1501427285778791
